$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Data")

# Update Avg_Time_ms values for the first two data rows.
$ws.Range("D2").Value = 11.1187151
$ws.Range("D3").Value = 45.296277
